# Add a new "2021" column (R) to the table on sheet1, mirroring the
# existing formatting of the adjacent "2020" column (P/Q) but giving the
# new data cell its own "0.0" number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year header in R4: copy Q4's format, then set the value.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# New data value in R5: copy Q5's format (font/border/alignment), then set
# its own value and a dedicated one-decimal number format.
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 102.20441221981518
$ws.Range("R5").NumberFormat = "0.0"

# Move the active selection, matching the saved view state after editing.
$ws.Range("S9").Select()
